# Apply data updates to the "earnings_debt" worksheet for South Korea / Bank (Money Center)
# per the commit "Updated capital structure database".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# --- Updated / newly populated cells ---
# Row 2
Set-Cell 2 4 0.1034
Set-Cell 2 5 0.1176
Set-Cell 2 6 0.003840000000000001
Set-Cell 2 9 0
Set-Cell 2 10 0
Set-Cell 2 11 10685.7
Set-Cell 2 12 0.2635730618104322
Set-Cell 2 13 3786.9389
Set-Cell 2 14 0.07083050717481408
Set-Cell 2 15 0.3543931515951224
Set-Cell 2 16 3005.8819
Set-Cell 2 17 0.05622169913662821
Set-Cell 2 18 0.2812994843575994
Set-Cell 2 19 781.0569999999999
Set-Cell 2 20 0.2062502249508171
Set-Cell 2 21 52434.3
Set-Cell 2 22 0.9807256363065044
Set-Cell 2 23 0.09089691771384896
Set-Cell 2 24 0.1801609094600408
Set-Cell 2 25 -0.0892639917461918
Set-Cell 2 26 0.07657320745180567
Set-Cell 2 27 0
Set-Cell 2 28 0.03615114117842917
Set-Cell 2 29 -0.03615114117842917
Set-Cell 2 30 513330.5
Set-Cell 2 31 0
Set-Cell 2 32 513330.5
Set-Cell 2 33 460896.2
Set-Cell 2 34 0.9056717654504192
Set-Cell 2 35 0.7737850916616307
Set-Cell 2 36 0.896055882930471
Set-Cell 2 37 0.7543709355455184
# Row 3
Set-Cell 3 4 0.194
Set-Cell 3 5 0.0535
Set-Cell 3 6 0.0345
Set-Cell 3 9 0
Set-Cell 3 10 0
Set-Cell 3 11 1221.9
Set-Cell 3 12 0.2256468024597884
Set-Cell 3 13 424.7525
Set-Cell 3 14 0.07078146611341632
Set-Cell 3 15 0.3476164170554055
Set-Cell 3 16 424.7525
Set-Cell 3 17 0.07078146611341632
Set-Cell 3 18 0.3476164170554055
Set-Cell 3 19 0
Set-Cell 3 20 0
Set-Cell 3 21 13552.8
Set-Cell 3 22 2.258461230815378
Set-Cell 3 23 0.06613909833447905
Set-Cell 3 24 0.5247966542314239
Set-Cell 3 25 -0.4586575558969448
Set-Cell 3 26 0.03745892737321961
Set-Cell 3 27 0
Set-Cell 3 28 0.03601525220688034
Set-Cell 3 29 -0.03601525220688034
Set-Cell 3 30 157470
Set-Cell 3 31 0
Set-Cell 3 32 157470
Set-Cell 3 33 143917.2
Set-Cell 3 34 0.9632907141271015
Set-Cell 3 35 0.8797231277182338
Set-Cell 3 36 0.959972144791056
Set-Cell 3 37 0.8698704779243104
# Row 4
Set-Cell 4 2 "Woori Financial Group Inc. (KOSE:A316140)"
Set-Cell 4 4 0.0868
Set-Cell 4 5 0.148
Set-Cell 4 6 -0.00185
Set-Cell 4 11 1157.2
Set-Cell 4 12 0.1949756533167091
Set-Cell 4 13 1160.5
Set-Cell 4 14 0.179688467731946
Set-Cell 4 15 1.002851711026616
Set-Cell 4 16 508.7
Set-Cell 4 17 0.07876563854824725
Set-Cell 4 18 0.4395955755271344
Set-Cell 4 19 651.8
Set-Cell 4 20 0.561654459284791
Set-Cell 4 21 7570.6
Set-Cell 4 22 1.172209835253314
Set-Cell 4 23 0.06654705219331539
Set-Cell 4 24 0.1757746341999422
Set-Cell 4 25 -0.1092275820066268
Set-Cell 4 26 0.1016458297653708
Set-Cell 4 28 0.036131768390056
Set-Cell 4 29 -0.036131768390056
Set-Cell 4 30 48900.1
Set-Cell 4 32 48900.1
Set-Cell 4 33 41329.5
Set-Cell 4 34 0.8833349892067162
Set-Cell 4 35 0.6847698533839325
Set-Cell 4 36 0.8648528183912664
Set-Cell 4 37 0.6473878609783147
# Row 5
Set-Cell 5 2 "JB Financial Group Co., Ltd. (KOSE:A175330)"
Set-Cell 5 4 0.131
Set-Cell 5 5 -0.103
Set-Cell 5 6 0.009000000000000001
Set-Cell 5 11 297.1
Set-Cell 5 12 0.2102767357916342
Set-Cell 5 13 50.1294
Set-Cell 5 14 0.04976610741586419
Set-Cell 5 15 0.1687290474587681
Set-Cell 5 16 50.1294
Set-Cell 5 17 0.04976610741586419
Set-Cell 5 18 0.1687290474587681
Set-Cell 5 19 0
Set-Cell 5 20 0
Set-Cell 5 21 2340.1
Set-Cell 5 22 2.32314107018763
Set-Cell 5 23 0.1033894766146993
Set-Cell 5 24 0.1845471847201393
Set-Cell 5 25 -0.08115770810543994
Set-Cell 5 26 0.1864254707147475
Set-Cell 5 28 0.03614483028868434
Set-Cell 5 29 -0.03614483028868434
Set-Cell 5 30 8099.5
Set-Cell 5 32 8099.5
Set-Cell 5 33 5759.4
Set-Cell 5 34 0.8893903456757588
Set-Cell 5 35 0.7041696371129001
Set-Cell 5 36 0.8511386643415549
Set-Cell 5 37 0.6286113445607449
# Row 6
Set-Cell 6 2 "Hana Financial Group Inc. (KOSE:A086790)"
Set-Cell 6 4 0.0225
Set-Cell 6 5 0.184
Set-Cell 6 6 0.08500000000000001
Set-Cell 6 11 2108.5
Set-Cell 6 12 0.2511733733590641
Set-Cell 6 13 555.157
Set-Cell 6 14 0.06005462884835897
Set-Cell 6 15 0.2632947593075646
Set-Cell 6 16 555.1
Set-Cell 6 17 0.06004846281993034
Set-Cell 6 18 0.2632677258714726
Set-Cell 6 19 0.05700000000001637
Set-Cell 6 20 0.0001026736580823377
Set-Cell 6 21 13300.6
Set-Cell 6 22 1.438804872244218
Set-Cell 6 23 0.08970775311541392
Set-Cell 6 24 0.193987624075191
Set-Cell 6 25 -0.1042798709597771
Set-Cell 6 26 0.1066739522708212
Set-Cell 6 28 0.03615745206817401
Set-Cell 6 29 -0.03615745206817401
Set-Cell 6 30 78998.89999999999
Set-Cell 6 32 78998.89999999999
Set-Cell 6 33 65698.29999999999
Set-Cell 6 34 0.8952416676204712
Set-Cell 6 35 0.7451953099206685
Set-Cell 6 36 0.8766494312306101
Set-Cell 6 37 0.7086400231257766
# Row 7
Set-Cell 7 4 0.07099999999999999
Set-Cell 7 5 0.0872
Set-Cell 7 6 -0.0222
Set-Cell 7 9 0
Set-Cell 7 10 0
Set-Cell 7 11 2969.9
Set-Cell 7 12 0.3020800488226619
Set-Cell 7 13 842.8
Set-Cell 7 14 0.05539780196666141
Set-Cell 7 15 0.283780598673356
Set-Cell 7 16 713.6
Set-Cell 7 17 0.04690540043119314
Set-Cell 7 18 0.240277450419206
Set-Cell 7 19 129.1999999999999
Set-Cell 7 20 0.153298528713811
Set-Cell 7 21 8290.4
Set-Cell 7 22 0.5449334805700162
Set-Cell 7 23 0.09226845120481179
Set-Cell 7 24 0.1662674643350106
Set-Cell 7 25 -0.07399901313019884
Set-Cell 7 26 0.07987329462412064
Set-Cell 7 27 0
Set-Cell 7 28 0.0373608292187612
Set-Cell 7 29 -0.0373608292187612
Set-Cell 7 30 107453.7
Set-Cell 7 31 0
Set-Cell 7 32 107453.7
Set-Cell 7 33 99163.3
Set-Cell 7 34 0.8759767272940709
Set-Cell 7 35 0.7314597173238685
Set-Cell 7 36 0.8669871276455298
Set-Cell 7 37 0.7153983725877932
# Row 8
Set-Cell 8 4 0.12
Set-Cell 8 5 0.17
Set-Cell 8 6 -0.00132
Set-Cell 8 9 0
Set-Cell 8 10 0
Set-Cell 8 11 2931.1
Set-Cell 8 12 0.3068411410625491
Set-Cell 8 13 753.6
Set-Cell 8 14 0.04849296028416257
Set-Cell 8 15 0.2571048411859029
Set-Cell 8 16 753.6
Set-Cell 8 17 0.04849296028416257
Set-Cell 8 18 0.2571048411859029
Set-Cell 8 19 0
Set-Cell 8 20 0
Set-Cell 8 21 7379.8
Set-Cell 8 22 0.4748783815088415
Set-Cell 8 23 0.092086082312284
Set-Cell 8 24 0.1694510377566659
Set-Cell 8 25 -0.07736495544438186
Set-Cell 8 26 0.08154939148488614
Set-Cell 8 27 0
Set-Cell 8 28 0.03737000828267786
Set-Cell 8 29 -0.03737000828267786
Set-Cell 8 30 112408.3
Set-Cell 8 31 0
Set-Cell 8 32 112408.3
Set-Cell 8 33 105028.5
Set-Cell 8 34 0.8785419468896519
Set-Cell 8 35 0.7565756307752779
Set-Cell 8 36 0.8711077234676604
Set-Cell 8 37 0.7438526636509856

# --- Cells removed entirely (columns AN/AP dropped for these rows) ---
# Row 2
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
# Row 3
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
# Row 7
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()
# Row 8
$ws.Range("AN8").ClearContents()
$ws.Range("AP8").ClearContents()

Write-Host "South Korea bank money center data updated."
